# Generate Report for Handback
# Adds two new handed-back files (3d6ec34f... and f7ea0717...) as new rows
# to the Overview sheet and the two per-locale detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"
$reason = "Include"

# ---- new entity identifiers (known from the target change) ----
$md1 = "3d6ec34f-ab4b-49a3-9e3f-97645cb1b407.md"
$md2 = "f7ea0717-8b67-4abc-bafb-a6581a171753.md"

$xlf1_zh = "3d6ec34f-ab4b-49a3-9e3f-97645cb1b407.09484f996a63f0d8d9cfce5af40b646cdb17ba91.zh-cn.xlf"
$xlf2_zh = "f7ea0717-8b67-4abc-bafb-a6581a171753.ab5ba992854bb62ef2223e6292a91efc655a2ba6.zh-cn.xlf"
$xlf1_de = "3d6ec34f-ab4b-49a3-9e3f-97645cb1b407.09484f996a63f0d8d9cfce5af40b646cdb17ba91.de-de.xlf"
$xlf2_de = "f7ea0717-8b67-4abc-bafb-a6581a171753.ab5ba992854bb62ef2223e6292a91efc655a2ba6.de-de.xlf"

$dtOff1_zh = "2016-02-17 06:42:29"
$dtBack1_zh = "2016-02-17 06:43:16"
$dtOff1_de = "2016-02-17 06:42:40"
$dtBack1_de = "2016-02-17 06:43:34"

# rows 6 and 7 reuse the same handoff/handback timestamp pair for both files
$dtOff2_zh = $dtOff1_zh
$dtBack2_zh = $dtBack1_zh
$dtOff2_de = $dtOff1_de
$dtBack2_de = $dtBack1_de

$base1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/09484f996a63f0d8d9cfce5af40b646cdb17ba91/e2e"
$base2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/ab5ba992854bb62ef2223e6292a91efc655a2ba6/e2e"

# ---------------------------------------------------------------
# Sheet "Overview": two new rows, columns A (link+name), B, C (status)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B6").Value = $status
$wsOverview.Range("C6").Value = $status
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "$base1/$md1", "", "", $md1) | Out-Null

$wsOverview.Range("B7").Value = $status
$wsOverview.Range("C7").Value = $status
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "$base2/$md2", "", "", $md2) | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn": two new rows across columns A..H
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$hoff1_zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09484f996a63f0d8d9cfce5af40b646cdb17ba91/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$xlf1_zh"
$hback1_zh = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/09484f996a63f0d8d9cfce5af40b646cdb17ba91/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$xlf1_zh"
$hoff2_zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ab5ba992854bb62ef2223e6292a91efc655a2ba6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$xlf2_zh"
$hback2_zh = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ab5ba992854bb62ef2223e6292a91efc655a2ba6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$xlf2_zh"

# Row 6 (3d6ec34f...)
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "$base1/$md1", "", "", $md1) | Out-Null
$wsZh.Range("B6").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), $hoff1_zh, "", "", $xlf1_zh) | Out-Null
$wsZh.Range("D6").Value = $dtOff1_zh
$wsZh.Hyperlinks.Add($wsZh.Range("E6"), "$base1/$md1", "", "", $md1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F6"), $hoff1_zh, "", "", $xlf1_zh) | Out-Null
$wsZh.Range("G6").Value = $dtBack1_zh
$wsZh.Range("H6").Value = $reason

# Row 7 (f7ea0717...)
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "$base2/$md2", "", "", $md2) | Out-Null
$wsZh.Range("B7").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), $hoff2_zh, "", "", $xlf2_zh) | Out-Null
$wsZh.Range("D7").Value = $dtOff2_zh
$wsZh.Hyperlinks.Add($wsZh.Range("E7"), "$base2/$md2", "", "", $md2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F7"), $hoff2_zh, "", "", $xlf2_zh) | Out-Null
$wsZh.Range("G7").Value = $dtBack2_zh
$wsZh.Range("H7").Value = $reason

# ---------------------------------------------------------------
# Sheet "de-de": two new rows across columns A..H
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$hoff1_de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09484f996a63f0d8d9cfce5af40b646cdb17ba91/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$xlf1_de"
$hback1_de = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/09484f996a63f0d8d9cfce5af40b646cdb17ba91/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$xlf1_de"
$hoff2_de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ab5ba992854bb62ef2223e6292a91efc655a2ba6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$xlf2_de"
$hback2_de = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ab5ba992854bb62ef2223e6292a91efc655a2ba6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$xlf2_de"

# Row 6 (3d6ec34f...)
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "$base1/$md1", "", "", $md1) | Out-Null
$wsDe.Range("B6").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), $hoff1_de, "", "", $xlf1_de) | Out-Null
$wsDe.Range("D6").Value = $dtOff1_de
$wsDe.Hyperlinks.Add($wsDe.Range("E6"), "$base1/$md1", "", "", $md1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F6"), $hoff1_de, "", "", $xlf1_de) | Out-Null
$wsDe.Range("G6").Value = $dtBack1_de
$wsDe.Range("H6").Value = $reason

# Row 7 (f7ea0717...)
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "$base2/$md2", "", "", $md2) | Out-Null
$wsDe.Range("B7").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), $hoff2_de, "", "", $xlf2_de) | Out-Null
$wsDe.Range("D7").Value = $dtOff2_de
$wsDe.Hyperlinks.Add($wsDe.Range("E7"), "$base2/$md2", "", "", $md2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F7"), $hoff2_de, "", "", $xlf2_de) | Out-Null
$wsDe.Range("G7").Value = $dtBack2_de
$wsDe.Range("H7").Value = $reason

Write-Output "Report rows added for Overview/zh-cn/de-de"
